$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 438.77777
$ws.Range("I2").Value = 368.625
$ws.Range("K2").Value = 368.625
$ws.Range("M2").Value = -255.625
$ws.Range("H17").Value = 1572.45
$ws.Range("J17").Value = 1608.8334
$ws.Range("L17").Value = 4826.5002
$ws.Range("N17").Value = -5162.5002
$ws.Range("H18").Value = 142858270
$ws.Range("I18").Value = 574.4
$ws.Range("J18").Value = 500002500
$ws.Range("K18").Value = 574.4
$ws.Range("L18").Value = 500002500
$ws.Range("M18").Value = -290.4
$ws.Range("N18").Value = -500003068
$ws.Range("H29").Value = 140
$ws.Range("I29").Value = 140
$ws.Range("K29").Value = 420
$ws.Range("M29").Value = -139
$ws.Range("H80").Value = 774.2727
$ws.Range("I80").Value = 439.66666
$ws.Range("K80").Value = 1318.99998
$ws.Range("M80").Value = -320.9999800000001
$ws.Range("H83").Value = 774.2727
$ws.Range("I83").Value = 439.66666
$ws.Range("K83").Value = 3956.99994
$ws.Range("M83").Value = 1035.00006
$ws.Range("H100").Value = 1873.5
$ws.Range("I100").Value = 1873.5
$ws.Range("K100").Value = 1873.5
$ws.Range("M100").Value = -1332.5
$ws.Range("H112").Value = 1076.7778
$ws.Range("J112").Value = 1538.2
$ws.Range("L112").Value = 4614.6
$ws.Range("N112").Value = -6830.6
$ws.Range("H132").Value = 6346.64
$ws.Range("I132").Value = 6346.64
$ws.Range("K132").Value = 19039.92
$ws.Range("M132").Value = -16509.92

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 77049
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 77049
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 77049
$ws.Range("M44").Value = ""
$ws.Range("N44").Value = -78025
$ws.Range("H97").Value = 732.86664
$ws.Range("I97").Value = 628.0714
$ws.Range("J97").Value = 2200
$ws.Range("K97").Value = 628.0714
$ws.Range("L97").Value = 2200
$ws.Range("M97").Value = -132.0714
$ws.Range("N97").Value = -3192

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 105269290
$ws.Range("I94").Value = 133340630
$ws.Range("K94").Value = 133340630
$ws.Range("M94").Value = -133340179
$ws.Range("H97").Value = 29999.5
$ws.Range("J97").Value = 29999.5
$ws.Range("L97").Value = 29999.5
$ws.Range("N97").Value = -31981.5
$ws.Range("H100").Value = 20000
$ws.Range("J100").Value = 20000
$ws.Range("L100").Value = 20000
$ws.Range("N100").Value = -22164

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2195
$ws.Range("J16").Value = 2148
$ws.Range("L16").Value = 2148
$ws.Range("N16").Value = -2722
$ws.Range("H31").Value = 5734.543
$ws.Range("I31").Value = 3732.65
$ws.Range("K31").Value = 3732.65
$ws.Range("M31").Value = -3437.65
$ws.Range("H34").Value = 5734.543
$ws.Range("I34").Value = 3732.65
$ws.Range("K34").Value = 3732.65
$ws.Range("M34").Value = -3530.65
$ws.Range("H58").Value = 2152.1177
$ws.Range("I58").Value = 1164.7
$ws.Range("J58").Value = 3562.7144
$ws.Range("K58").Value = 1164.7
$ws.Range("L58").Value = 3562.7144
$ws.Range("M58").Value = -961.7
$ws.Range("N58").Value = -3968.7144
$ws.Range("H94").Value = 2756
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 2756
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 2756
$ws.Range("M94").Value = ""
$ws.Range("N94").Value = -3658
$ws.Range("H99").Value = 2966.3333
$ws.Range("I99").Value = 1999.6666
$ws.Range("J99").Value = 3933
$ws.Range("K99").Value = 1999.6666
$ws.Range("L99").Value = 3933
$ws.Range("M99").Value = -501.6666
$ws.Range("N99").Value = -6929
$ws.Range("H113").Value = 2195
$ws.Range("J113").Value = 2148
$ws.Range("L113").Value = 2148
$ws.Range("N113").Value = -6488
$ws.Range("H122").Value = 2816.5
$ws.Range("I122").Value = 3088.6667
$ws.Range("K122").Value = 9266.000100000001
$ws.Range("M122").Value = -6816.000100000001
$ws.Range("H126").Value = 2966.3333
$ws.Range("I126").Value = 1999.6666
$ws.Range("J126").Value = 3933
$ws.Range("K126").Value = 5998.9998
$ws.Range("L126").Value = 11799
$ws.Range("M126").Value = -3528.9998
$ws.Range("N126").Value = -16739
$ws.Range("H132").Value = 19615186
$ws.Range("I132").Value = 5429.143
$ws.Range("K132").Value = 16287.429
$ws.Range("M132").Value = -13757.429
$ws.Range("H134").Value = 4197.0625
$ws.Range("I134").Value = 3575.2144
$ws.Range("K134").Value = 10725.6432
$ws.Range("M134").Value = -8190.643199999999
$ws.Range("H136").Value = 2152.1177
$ws.Range("I136").Value = 1164.7
$ws.Range("J136").Value = 3562.7144
$ws.Range("K136").Value = 3494.1
$ws.Range("L136").Value = 10688.1432
$ws.Range("M136").Value = -944.1000000000004
$ws.Range("N136").Value = -15788.1432
$ws.Range("H138").Value = 61999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 452
$ws.Range("J12").Value = 493
$ws.Range("L12").Value = 1479
$ws.Range("N12").Value = -1825
$ws.Range("H44").Value = 2036.0834
$ws.Range("J44").Value = 4490
$ws.Range("L44").Value = 13470
$ws.Range("N44").Value = -14266
$ws.Range("H80").Value = 19724.5
$ws.Range("I80").Value = 19700
$ws.Range("J80").Value = 19732.666
$ws.Range("K80").Value = 59100
$ws.Range("L80").Value = 59197.99800000001
$ws.Range("N80").Value = -61069.99800000001
$ws.Range("M80").Value = -58164
$ws.Range("H83").Value = 19724.5
$ws.Range("I83").Value = 19700
$ws.Range("J83").Value = 19732.666
$ws.Range("K83").Value = 177300
$ws.Range("L83").Value = 177593.994
$ws.Range("N83").Value = -186953.994
$ws.Range("M83").Value = -172620
$ws.Range("H114").Value = 3798.6
$ws.Range("I114").Value = 3249
$ws.Range("J114").Value = 4165
$ws.Range("K114").Value = 9747
$ws.Range("L114").Value = 12495
$ws.Range("M114").Value = -6493
$ws.Range("N114").Value = -19003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2067.1
$ws.Range("I7").Value = 1968.7142
$ws.Range("J7").Value = 2296.6667
$ws.Range("K7").Value = 1968.7142
$ws.Range("L7").Value = 2296.6667
$ws.Range("M7").Value = -1856.7142
$ws.Range("N7").Value = -2520.6667
$ws.Range("H36").Value = 74990
$ws.Range("J36").Value = 74990
$ws.Range("L36").Value = 74990
$ws.Range("N36").Value = -76114
$ws.Range("H55").Value = 725.1818
$ws.Range("I55").Value = 622
$ws.Range("J55").Value = 849
$ws.Range("K55").Value = 622
$ws.Range("L55").Value = 849
$ws.Range("M55").Value = -449
$ws.Range("N55").Value = -1195
$ws.Range("H122").Value = 13990.143
$ws.Range("I122").Value = 12396.8
$ws.Range("K122").Value = 37190.39999999999
$ws.Range("M122").Value = -34740.39999999999
$ws.Range("H126").Value = 2067.1
$ws.Range("I126").Value = 1968.7142
$ws.Range("J126").Value = 2296.6667
$ws.Range("K126").Value = 5906.142599999999
$ws.Range("L126").Value = 6890.000100000001
$ws.Range("M126").Value = -3436.142599999999
$ws.Range("N126").Value = -11830.0001
